$wb = $excel.ActiveWorkbook

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/38c5c55f9a3f0b4745e174e4c42c48d58e79ccc3/e2e/ba056d89-61db-4787-85ec-ff51c52bd823.md"
$hyperlinkDisplay = "ba056d89-61db-4787-85ec-ff51c52bd823.md"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("J2").Value = "ba056d89-61db-4787-85ec-ff51c52bd823.346578315c49711ff87c52feab484250854e17e6.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-20 05:01:10"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $hyperlinkTarget, "", "", $hyperlinkDisplay)
$wsZhCn.Range("I2").Font.Name = "Calibri"
$wsZhCn.Range("I2").Font.Size = 11
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(9).ColumnWidth = 40
$wsZhCn.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("J2").Value = "ba056d89-61db-4787-85ec-ff51c52bd823.346578315c49711ff87c52feab484250854e17e6.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-20 05:01:16"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $hyperlinkTarget, "", "", $hyperlinkDisplay)
$wsDeDe.Range("I2").Font.Name = "Calibri"
$wsDeDe.Range("I2").Font.Size = 11
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(9).ColumnWidth = 40
$wsDeDe.Columns.Item(10).ColumnWidth = 40
